# "update1.1" — turns the single login test case into a 2-row test suite
# (adds an "Action/Steps" + "Test result/Factual result" pair of columns,
#  rewrites the existing login case, and appends a "send comment" case).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft    = -4131
$xlVTop    = -4160
$greenFill = 5296274   # RGB(146, 208, 80) == 0x92D050, stored BGR-order as 5296274

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Pre-condition(Предусловие)"
$ws.Range("D1").Value = "Action(действия)/Steps"
$ws.Range("E1").Value = "Test data(тестовые данные)"
$ws.Range("F1").Value = "Expected result(ожидаемый результат)"
$ws.Range("G1").Value = "Test result/Factual result"

# ---------------------------------------------------------------------------
# Row 2 — existing "login" test case, reworked (values first)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1

$ws.Range("B2").Value = 'Открыть страницу "Вход в систему"'

$ws.Range("C2").Value = @'
1. Открыть  сайт Ситилинк
https://www.citilink.ru/
2. Пользователь должен быть уже зарегистрирован в системе

'@

$ws.Range("D2").Value = @'
1. Нажать на кнопку "Войти" 
в правой части сайта
2.Ввести Email и пароль
в поля ввода
3. Нажать на кнопку "Войти"
'@

$ws.Range("E2").Value = @'
"E-mail"- anglenich@mail.ru
"Password"-angela1919
'@

$ws.Range("F2").Value = @'
1. Окно "Вход в систему" открыто            2. Название окна "Вход в систему"           3. Логотип компании отображается в правом верхнем углу                                            4. На форме 2 поля "Email" и "Пароль"     5. Кнопка Вход доступна
6. Пользователя впустило в систему                                
7. Ссылка "забыть пароль" доступна       
'@

$ws.Range("G2").Value = @'
Пользователя впустило в 
систему 
Тест прошел - passed
'@

$ws.Rows(2).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 3 stays empty/unchanged (kept as a spacer row, same as before)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 4 — new "send comment" test case (values first)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 2

$ws.Range("B4").Value = "Тест отправки комментария"

$ws.Range("C4").Value = @'
1. Открыть  сайт Ситилинк
https://www.citilink.ru/
2. Перейти по сыллке 
с названием "Задать впорос" 
внизу страницы
'@

$ws.Range("D4").Value = @'
1. Заполнить форму 
отправки комментрая
2. Нажать на кнопку "Отправить"
'@

$ws.Range("E4").Value = @'
"Тип обращения"-Комментарий
"Контактное лицо"- Ангелина
"E-mail"- anglenich@mail.ru
"Сообщение"- Добрый день, уважаемые пользователи данного сайта!
Хочу поздравить всех с предстоящими 
праздниками и пожелать всем добра и
счастья! Уважайте друг друга)
'@

$ws.Range("F4").Value = @'
1. Сайт Про тестинг открыт и доступен
2. Страница "Вопросы, пожелания и заявки" открыта и доступна
3. Данные успешно введены
4. Комментарий успешно отправлен
5. Страница "Ваш запрос успешно отправлен!" открыта
'@

$ws.Range("G4").Value = @'
Запрос успешно отправлен
Тест прошел - passed
'@

$ws.Rows(4).RowHeight = 135

# ---------------------------------------------------------------------------
# Styles — applied in the same order the original authoring tool created
# them, so the new cellXfs entries line up at indices 3 / 4 / 5:
#   3 = horizontal left + vertical top   (B4)
#   4 = vertical top + wrap text         (C2,D2,E2,C4,D4,F4)
#   5 = fill(green) + vertical top + wrap text (G2,G4)
# ---------------------------------------------------------------------------
$ws.Range("A2").VerticalAlignment = $xlVTop
$ws.Range("B2").VerticalAlignment = $xlVTop
$ws.Range("A4").VerticalAlignment = $xlVTop

$ws.Range("B4").HorizontalAlignment = $xlLeft
$ws.Range("B4").VerticalAlignment = $xlVTop

$ws.Range("C2:E2").VerticalAlignment = $xlVTop
$ws.Range("C2:E2").WrapText = $true

$ws.Range("F2").WrapText = $true

$ws.Range("C4:F4").VerticalAlignment = $xlVTop
$ws.Range("C4:F4").WrapText = $true

$ws.Range("G2").VerticalAlignment = $xlVTop
$ws.Range("G2").WrapText = $true
$ws.Range("G2").Interior.Color = $greenFill

$ws.Range("G4").VerticalAlignment = $xlVTop
$ws.Range("G4").WrapText = $true
$ws.Range("G4").Interior.Color = $greenFill

# ---------------------------------------------------------------------------
# Column widths / page setup
# ---------------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 21.7
$ws.Columns("G:G").ColumnWidth = 22.6

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("J4").Select()
